$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("y_fitted_on_begin_2016")
$ws1.Range("B2").Value = 0.236171715240392
$ws1.Range("B3").Value = 72.04520429109607
$ws1.Range("B4").Value = 72.71598109934938
$ws1.Range("B5").Value = 72.83095061166003
$ws1.Range("B6").Value = 72.90085633178259
$ws1.Range("B7").Value = 74.12022053551081
$ws1.Range("B8").Value = 74.57755161591204
$ws1.Range("B9").Value = 74.08861328985996
$ws1.Range("B10").Value = 73.85658676077362
$ws1.Range("B11").Value = 74.10893133356851
$ws1.Range("B12").Value = 74.31266212202878
$ws1.Range("B13").Value = 74.95003124256243
$ws1.Range("B14").Value = 75.71827256281534
$ws1.Range("B15").Value = 76.8071028210129
$ws1.Range("B16").Value = 76.17424774185801
$ws1.Range("B17").Value = 75.8246590609198
$ws1.Range("B18").Value = 76.08090088275939
$ws1.Range("B19").Value = 76.0739359457107
$ws1.Range("B20").Value = 76.47310204024078
$ws1.Range("B21").Value = 77.28092901672069

$ws3 = $wb.Worksheets.Item("y_fitted_on_begin_2021")
$ws3.Range("B2").Value = 0.2333812346860189
$ws3.Range("B3").Value = 72.04240584069341
$ws3.Range("B4").Value = 72.62714614997304
$ws3.Range("B5").Value = 72.79546284961691
$ws3.Range("B6").Value = 72.92402834652951
$ws3.Range("B7").Value = 73.83778953530492
$ws3.Range("B8").Value = 74.16826759059731
$ws3.Range("B9").Value = 73.87879086581108
$ws3.Range("B10").Value = 73.86238933434508
$ws3.Range("B11").Value = 74.23131447621341
$ws3.Range("B12").Value = 74.60330339757037
$ws3.Range("B13").Value = 75.26236734201363
$ws3.Range("B14").Value = 75.9153930800461
$ws3.Range("B15").Value = 76.5398785482675
$ws3.Range("B16").Value = 75.96177511307265
$ws3.Range("B17").Value = 75.73743312717671
$ws3.Range("B18").Value = 75.98149298458719
$ws3.Range("B19").Value = 76.12579436919259
$ws3.Range("B20").Value = 76.61249940402857
$ws3.Range("B21").Value = 77.2728256104209
$ws3.Range("B22").Value = 77.50316998349699
$ws3.Range("B23").Value = 77.39142453335245
$ws3.Range("B24").Value = 77.6680658824814
$ws3.Range("B25").Value = 78.27328006136199
$ws3.Range("B26").Value = 78.87931812479869

